# Delete the last 4 sending-cluster rows (Resolving-Mac as a sending cluster no longer present)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("14:17").Delete()

# Update remaining rows 2-13 with the refreshed TPM-derived values
# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 7).Value = 4.343956666666666
$ws.Cells.Item(2, 8).Value = 13.03187
$ws.Cells.Item(2, 9).Value = 0.2551833209483726
$ws.Cells.Item(2, 10).Value = 0.2551833209483726
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 2.613621
$ws.Cells.Item(2, 14).Value = 7.840863
$ws.Cells.Item(2, 15).Value = 0.4199754212870037
$ws.Cells.Item(2, 16).Value = 0.4199754212870037
$ws.Cells.Item(2, 17).Value = 11.35345636709
$ws.Cells.Item(2, 18).Value = 102.18110730381
$ws.Cells.Item(2, 19).Value = 0.1071707227207095
$ws.Cells.Item(2, 20).Value = 0.1071707227207095

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 7).Value = 4.343956666666666
$ws.Cells.Item(3, 8).Value = 13.03187
$ws.Cells.Item(3, 9).Value = 0.2551833209483726
$ws.Cells.Item(3, 10).Value = 0.2551833209483726
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 3.430056
$ws.Cells.Item(3, 14).Value = 10.290168
$ws.Cells.Item(3, 15).Value = 0.551166069463788
$ws.Cells.Item(3, 16).Value = 0.5511660694637879
$ws.Cells.Item(3, 17).Value = 14.90001462824
$ws.Cells.Item(3, 18).Value = 134.10013165416
$ws.Cells.Item(3, 19).Value = 0.1406483879998309
$ws.Cells.Item(3, 20).Value = 0.1406483879998308

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 7).Value = 4.343956666666666
$ws.Cells.Item(4, 8).Value = 13.03187
$ws.Cells.Item(4, 9).Value = 0.2551833209483726
$ws.Cells.Item(4, 10).Value = 0.2551833209483726
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.1795943333333333
$ws.Cells.Item(4, 14).Value = 0.538783
$ws.Cells.Item(4, 15).Value = 0.02885850924920838
$ws.Cells.Item(4, 16).Value = 0.02885850924920838
$ws.Cells.Item(4, 17).Value = 0.7801500015788888
$ws.Cells.Item(4, 18).Value = 7.02135001421
$ws.Cells.Item(4, 19).Value = 0.007364210227832322
$ws.Cells.Item(4, 20).Value = 0.007364210227832322

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 7).Value = 6.440526333333334
$ws.Cells.Item(5, 8).Value = 19.321579
$ws.Cells.Item(5, 9).Value = 0.3783451411951115
$ws.Cells.Item(5, 10).Value = 0.3783451411951115
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 2.613621
$ws.Cells.Item(5, 14).Value = 7.840863
$ws.Cells.Item(5, 15).Value = 0.4199754212870037
$ws.Cells.Item(5, 16).Value = 0.4199754212870037
$ws.Cells.Item(5, 17).Value = 16.833094875853
$ws.Cells.Item(5, 18).Value = 151.497853882677
$ws.Cells.Item(5, 19).Value = 0.1588956600653078
$ws.Cells.Item(5, 20).Value = 0.1588956600653078

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 7).Value = 6.440526333333334
$ws.Cells.Item(6, 8).Value = 19.321579
$ws.Cells.Item(6, 9).Value = 0.3783451411951115
$ws.Cells.Item(6, 10).Value = 0.3783451411951115
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 3.430056
$ws.Cells.Item(6, 14).Value = 10.290168
$ws.Cells.Item(6, 15).Value = 0.551166069463788
$ws.Cells.Item(6, 16).Value = 0.5511660694637879
$ws.Cells.Item(6, 17).Value = 22.091365992808
$ws.Cells.Item(6, 18).Value = 198.822293935272
$ws.Cells.Item(6, 19).Value = 0.2085310043732315
$ws.Cells.Item(6, 20).Value = 0.2085310043732314

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 4).Value = "MuSCs"
$ws.Cells.Item(7, 7).Value = 6.440526333333334
$ws.Cells.Item(7, 8).Value = 19.321579
$ws.Cells.Item(7, 9).Value = 0.3783451411951115
$ws.Cells.Item(7, 10).Value = 0.3783451411951115
$ws.Cells.Item(7, 11).Value = 1
$ws.Cells.Item(7, 12).Value = 0.3333333333333333
$ws.Cells.Item(7, 13).Value = 0.1795943333333333
$ws.Cells.Item(7, 14).Value = 0.538783
$ws.Cells.Item(7, 15).Value = 0.02885850924920838
$ws.Cells.Item(7, 16).Value = 0.02885850924920838
$ws.Cells.Item(7, 17).Value = 1.156682033150778
$ws.Cells.Item(7, 18).Value = 10.410138298357
$ws.Cells.Item(7, 19).Value = 0.01091847675657217
$ws.Cells.Item(7, 20).Value = 0.01091847675657217

# Row 8
$ws.Cells.Item(8, 1).Value = "MuSCs"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 7).Value = 0.3495363333333334
$ws.Cells.Item(8, 8).Value = 1.048609
$ws.Cells.Item(8, 9).Value = 0.02053331770470026
$ws.Cells.Item(8, 10).Value = 0.02053331770470026
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 2.613621
$ws.Cells.Item(8, 14).Value = 7.840863
$ws.Cells.Item(8, 15).Value = 0.4199754212870037
$ws.Cells.Item(8, 16).Value = 0.4199754212870037
$ws.Cells.Item(8, 17).Value = 0.9135555010630001
$ws.Cells.Item(8, 18).Value = 8.221999509567
$ws.Cells.Item(8, 19).Value = 0.008623488753451383
$ws.Cells.Item(8, 20).Value = 0.008623488753451383

# Row 9
$ws.Cells.Item(9, 1).Value = "MuSCs"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 7).Value = 0.3495363333333334
$ws.Cells.Item(9, 8).Value = 1.048609
$ws.Cells.Item(9, 9).Value = 0.02053331770470026
$ws.Cells.Item(9, 10).Value = 0.02053331770470026
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 3.430056
$ws.Cells.Item(9, 14).Value = 10.290168
$ws.Cells.Item(9, 15).Value = 0.551166069463788
$ws.Cells.Item(9, 16).Value = 0.5511660694637879
$ws.Cells.Item(9, 17).Value = 1.198929197368
$ws.Cells.Item(9, 18).Value = 10.790362776312
$ws.Cells.Item(9, 19).Value = 0.01131726801235085
$ws.Cells.Item(9, 20).Value = 0.01131726801235085

# Row 10
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 4).Value = "MuSCs"
$ws.Cells.Item(10, 7).Value = 0.3495363333333334
$ws.Cells.Item(10, 8).Value = 1.048609
$ws.Cells.Item(10, 9).Value = 0.02053331770470026
$ws.Cells.Item(10, 10).Value = 0.02053331770470026
$ws.Cells.Item(10, 11).Value = 1
$ws.Cells.Item(10, 12).Value = 0.3333333333333333
$ws.Cells.Item(10, 13).Value = 0.1795943333333333
$ws.Cells.Item(10, 14).Value = 0.538783
$ws.Cells.Item(10, 15).Value = 0.02885850924920838
$ws.Cells.Item(10, 16).Value = 0.02885850924920838
$ws.Cells.Item(10, 17).Value = 0.06277474476077778
$ws.Cells.Item(10, 18).Value = 0.5649727028470001
$ws.Cells.Item(10, 19).Value = 0.0005925609388980265
$ws.Cells.Item(10, 20).Value = 0.0005925609388980265

# Row 11
$ws.Cells.Item(11, 1).Value = "Resolving-Mac"
$ws.Cells.Item(11, 4).Value = "ECs"
$ws.Cells.Item(11, 7).Value = 5.888867
$ws.Cells.Item(11, 8).Value = 17.666601
$ws.Cells.Item(11, 9).Value = 0.3459382201518156
$ws.Cells.Item(11, 10).Value = 0.3459382201518156
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 2.613621
$ws.Cells.Item(11, 14).Value = 7.840863
$ws.Cells.Item(11, 15).Value = 0.4199754212870037
$ws.Cells.Item(11, 16).Value = 0.4199754212870037
$ws.Cells.Item(11, 17).Value = 15.391266457407
$ws.Cells.Item(11, 18).Value = 138.521398116663
$ws.Cells.Item(11, 19).Value = 0.145285549747535
$ws.Cells.Item(11, 20).Value = 0.145285549747535

# Row 12
$ws.Cells.Item(12, 1).Value = "Resolving-Mac"
$ws.Cells.Item(12, 4).Value = "FAPs"
$ws.Cells.Item(12, 7).Value = 5.888867
$ws.Cells.Item(12, 8).Value = 17.666601
$ws.Cells.Item(12, 9).Value = 0.3459382201518156
$ws.Cells.Item(12, 10).Value = 0.3459382201518156
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 3.430056
$ws.Cells.Item(12, 14).Value = 10.290168
$ws.Cells.Item(12, 15).Value = 0.551166069463788
$ws.Cells.Item(12, 16).Value = 0.5511660694637879
$ws.Cells.Item(12, 17).Value = 20.199143586552
$ws.Cells.Item(12, 18).Value = 181.792292278968
$ws.Cells.Item(12, 19).Value = 0.1906694090783748
$ws.Cells.Item(12, 20).Value = 0.1906694090783748

# Row 13
$ws.Cells.Item(13, 1).Value = "Resolving-Mac"
$ws.Cells.Item(13, 4).Value = "MuSCs"
$ws.Cells.Item(13, 7).Value = 5.888867
$ws.Cells.Item(13, 8).Value = 17.666601
$ws.Cells.Item(13, 9).Value = 0.3459382201518156
$ws.Cells.Item(13, 10).Value = 0.3459382201518156
$ws.Cells.Item(13, 11).Value = 1
$ws.Cells.Item(13, 12).Value = 0.3333333333333333
$ws.Cells.Item(13, 13).Value = 0.1795943333333333
$ws.Cells.Item(13, 14).Value = 0.538783
$ws.Cells.Item(13, 15).Value = 0.02885850924920838
$ws.Cells.Item(13, 16).Value = 0.02885850924920838
$ws.Cells.Item(13, 17).Value = 1.057607142953667
$ws.Cells.Item(13, 18).Value = 9.518464286583001
$ws.Cells.Item(13, 19).Value = 0.009983261325905855
$ws.Cells.Item(13, 20).Value = 0.009983261325905855
